$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: put literal text that LOOKS like a percentage (e.g. "88.7%") into a
# cell without Excel re-interpreting it as a numeric percentage. We briefly
# mark the cell as Text, assign the literal string, then paste the original
# (General) number formatting back on top so the cell's look/format is
# unchanged from before - only its text content differs.
function Set-LiteralText($rangeAddr, $text, $formatSourceAddr) {
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value2 = $text
    $ws.Range($formatSourceAddr).Copy() | Out-Null
    $c.PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# 1) Swap "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com" for
#    every "Recorded By" (column G) cell that currently has that exact text.
# ---------------------------------------------------------------------------
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$swapped = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
        $swapped = $swapped + 1
    }
}
Write-Host "Swapped Recorded-By order in $swapped cells"

# ---------------------------------------------------------------------------
# 2) Class Statistics block (K/L columns) numeric refresh.
# ---------------------------------------------------------------------------
$ws.Range("L6").Value2 = 141      # Recorded Sessions
$ws.Range("L7").Value2 = 0        # Missing Sessions
Set-LiteralText "L9" "88.7%" "L5"   # Coverage %
Set-LiteralText "L10" "71.4%" "L5"  # Average Attendance %

# ---------------------------------------------------------------------------
# 3) Per-group summary rows 18-20 (columns O/P/R/S).
# ---------------------------------------------------------------------------
$ws.Range("O18").Value2 = 23
$ws.Range("P18").Value2 = 0
Set-LiteralText "R18" "88.5%" "Q18"

$ws.Range("O19").Value2 = 23
$ws.Range("P19").Value2 = 0
Set-LiteralText "R19" "88.5%" "Q19"
Set-LiteralText "S19" "73.2%" "Q19"

$ws.Range("O20").Value2 = 23
$ws.Range("P20").Value2 = 0
Set-LiteralText "R20" "88.5%" "Q20"
Set-LiteralText "S20" "81.1%" "Q20"

# ---------------------------------------------------------------------------
# 4) Rows 105 / 131 / 157 were newly-recorded sessions: restyle them like the
#    rest of the data rows (formerly a distinct "not recorded" pink style,
#    now the standard style used by e.g. row 2) and fill in the now-known
#    Recorded-By / Students / Status values.
# ---------------------------------------------------------------------------
$formatSource = $ws.Range("A2:I2")

$formatSource.Copy() | Out-Null
$ws.Range("A105:I105").PasteSpecial(-4122) | Out-Null
$formatSource.Copy() | Out-Null
$ws.Range("A131:I131").PasteSpecial(-4122) | Out-Null
$formatSource.Copy() | Out-Null
$ws.Range("A157:I157").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("G105").Value2 = "dnasr281@gmail.com"
$ws.Range("H105").Value2 = "43/56"
$ws.Range("I105").Value2 = "Recorded"

$ws.Range("G131").Value2 = "dnasr281@gmail.com"
$ws.Range("H131").Value2 = "44/56"
$ws.Range("I131").Value2 = "Recorded"

$ws.Range("G157").Value2 = "dnasr281@gmail.com"
$ws.Range("H157").Value2 = "47/57"
$ws.Range("I157").Value2 = "Recorded"

# ---------------------------------------------------------------------------
# 5) Column I ("Status") width narrows from 14 to 10 (stored OOXML "width" is
#    ColumnWidth + ~0.83 padding for this workbook's default font, so ask
#    for 9.17 to land on a stored value of exactly 10).
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 9.17

Write-Host "Done"
